# "fixed the vignette to use the new syntax"
# Prepend a "." run to several function-name text boxes (new syntax uses a
# leading dot, e.g. ".vztdraw_...") and grow the textboxes' width to fit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- TextBox 17 (id=18) : "vztdraw_cumulative_intensity" -------------------
$sh18 = $s.Shapes.Item(10)
$sh18.Left = 48.64834976196289
$sh18.Width = 265.621826171875
$sh18.TextFrame.TextRange.InsertBefore(".")

# --- TextBox 18 (id=19) : "vztdraw_intensity" -------------------------------
$sh19 = $s.Shapes.Item(11)
$sh19.Left = 369.83197021484375
$sh19.Width = 173.08685302734375
$sh19.TextFrame.TextRange.InsertBefore(".")

# --- TextBox 19 (id=20) : already ".vztdraw_sc_step_regular_cpp" -----------
# only the width grows; text already carries the leading dot.
$sh20 = $s.Shapes.Item(12)
$sh20.Width = 261.7297668457031

# --- TextBox 24 (id=25) : two paragraphs ------------------------------------
$sh25 = $s.Shapes.Item(15)
$sh25.Width = 264.32427978515625
$tr25 = $sh25.TextFrame.TextRange
$tr25.Paragraphs(1, 1).InsertBefore(".")
$tr25.Paragraphs(2, 1).InsertBefore(".")

# --- TextBox 26 (id=27) : "vztdraw_intensity_step_regular" -----------------
$sh27 = $s.Shapes.Item(17)
$sh27.TextFrame.TextRange.InsertBefore(".")
